{"js": "// Word \"Use case 5\" update: transpose (\"transponeren\") -> inverse\n// (\"inverteren\"/\"inverse\") rewrite, per the supplied diff, plus a new\n// \"Gebruiker klikt op voltooien\" bullet and the removal of the\n// \"Gebruiker bekijkt de animatie...\" bullet.\n//\n// Every edit below is scoped to the paragraphs that sit between the\n// \"Use case 5\" heading and the next \"Use case\" heading, so the\n// near-duplicate text shared with the neighbouring use cases (4 and 6)\n// is left untouched.\n\nfunction normalize(t) {\n  return t.replace(/\\s+/g, \" \").trim();\n}\n\n// Re-(re)loads context.document.body.paragraphs text and returns\n// { items, startIdx, endIdx } describing the live paragraph list and the\n// [startIdx, endIdx) window that belongs to the \"Use case 5\" section.\nasync function loadUseCase5Window() {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n\n  const items = paragraphs.items;\n  let headingIdxs = [];\n  for (let i = 0; i < items.length; i++) {\n    if (/^Use case\\s*\\d+$/.test(normalize(items[i].text))) {\n      headingIdxs.push(i);\n    }\n  }\n\n  let startIdx = -1;\n  let endIdx = items.length;\n  for (let k = 0; k < headingIdxs.length; k++) {\n    const idx = headingIdxs[k];\n    if (normalize(items[idx].text) === \"Use case5\") {\n      startIdx = idx;\n      endIdx = k + 1 < headingIdxs.length ? headingIdxs[k + 1] : items.length;\n      break;\n    }\n  }\n\n  if (startIdx === -1) {\n    throw new Error(\"Could not locate the 'Use case 5' heading paragraph.\");\n  }\n\n  return { items, startIdx, endIdx };\n}\n\n// Finds, within the Use case 5 window, the first paragraph whose\n// (whitespace-normalized) text equals `expectedText` and replaces its\n// whole text with `newText`. Re-loads the window fresh each call so\n// earlier edits in this script are always visible to later searches.\nasync function replaceParagraphText(expectedText, newText) {\n  const { items, startIdx, endIdx } = await loadUseCase5Window();\n  for (let i = startIdx; i < endIdx; i++) {\n    if (normalize(items[i].text) === normalize(expectedText)) {\n      items[i].insertText(newText, Word.InsertLocation.replace);\n      await context.sync();\n      return;\n    }\n  }\n  throw new Error(\n    `Paragraph with text \"${expectedText}\" not found in Use case 5 section.`\n  );\n}\n\n// 1) \"Adjunct uitleg\" -> \"Inverse uitleg\"\nawait replaceParagraphText(\"Adjunct uitleg\", \"Inverse uitleg\");\n\n// 2) First \"Gebruiker weet hoe hij matrices moet transponeren\"\n//    (the \"Doelstelling\" bullet) -> \"...inverteren\"\nawait replaceParagraphText(\n  \"Gebruiker weet hoe hij matrices moet transponeren\",\n  \"Gebruiker weet hoe hij matrices moet inverteren\"\n);\n\n// 3) \"Basiskennis matrices\" -> append \"+determinanten kunnen uitrekenen\"\nawait replaceParagraphText(\n  \"Basiskennis matrices\",\n  \"Basiskennis matrices+determinanten kunnen uitrekenen\"\n);\n\n// 4) Second \"Gebruiker weet hoe hij matrices moet transponeren\"\n//    (the \"Postcondities\" bullet) -> \"...inverteren\". Because step 2\n//    above already turned the first occurrence into \"...inverteren\",\n//    this search (re-loaded fresh) now lands on the second bullet.\nawait replaceParagraphText(\n  \"Gebruiker weet hoe hij matrices moet transponeren\",\n  \"Gebruiker weet hoe hij matrices moet inverteren\"\n);\n\n// 5) \"Algemene definitie over transponeren van matrices\"\n//    -> \"Algemene definitie over inverteren van matrices\"\nawait replaceParagraphText(\n  \"Algemene definitie over transponeren van matrices\",\n  \"Algemene definitie over inverteren van matrices\"\n);\n\n// 6) \"Stap per stap animatie tonen\"\n//    -> \"Stap per stap bewerkingen uitleggen (adjunct berekenen)\"\nawait replaceParagraphText(\n  \"Stap per stap animatie tonen\",\n  \"Stap per stap bewerkingen uitleggen (adjunct berekenen)\"\n);\n\n// 7) Delete \"Gebruiker bekijkt de animatie en de extra uitleg stap per stap\"\n{\n  const { items, startIdx, endIdx } = await loadUseCase5Window();\n  let found = false;\n  for (let i = startIdx; i < endIdx; i++) {\n    if (\n      normalize(items[i].text) ===\n      \"Gebruiker bekijkt de animatie en de extra uitleg stap per stap\"\n    ) {\n      items[i].delete();\n      found = true;\n      break;\n    }\n  }\n  if (!found) {\n    throw new Error(\n      \"Paragraph 'Gebruiker bekijkt de animatie...' not found to delete.\"\n    );\n  }\n  await context.sync();\n}\n\n// 8) Insert new paragraph \"Gebruiker klikt op voltooien\" right after\n//    \"Gebruiker klikt op next knop om naar volgende stap te gaan\" and\n//    before \"Gebruiker krijgt matrix\".\n{\n  const { items, startIdx, endIdx } = await loadUseCase5Window();\n  let found = false;\n  for (let i = startIdx; i < endIdx; i++) {\n    if (\n      normalize(items[i].text) ===\n      \"Gebruiker klikt op next knop om naar volgende stap te gaan\"\n    ) {\n      items[i].insertParagraph(\n        \"Gebruiker klikt op voltooien\",\n        Word.InsertLocation.after\n      );\n      found = true;\n      break;\n    }\n  }\n  if (!found) {\n    throw new Error(\n      \"Paragraph 'Gebruiker klikt op next knop...' not found to insert after.\"\n    );\n  }\n  await context.sync();\n}\n\n// 9) \"Geruiker stelt zelf de getransponeerde \" -> \"Geruiker stelt zelf de inverse op \"\nawait replaceParagraphText(\n  \"Geruiker stelt zelf de getransponeerde\",\n  \"Geruiker stelt zelf de inverse op \"\n);\n", "ps1": "# Word \"Use case 5\" update: transpose (\"transponeren\") -> inverse\n# (\"inverteren\"/\"inverse\") rewrite, per the supplied diff, plus a new\n# \"Gebruiker klikt op voltooien\" bullet and the removal of the\n# \"Gebruiker bekijkt de animatie...\" bullet.\n#\n# Every edit below is scoped to the paragraphs that sit between the\n# \"Use case 5\" heading and the next \"Use case\" heading, so the\n# near-duplicate text shared with the neighbouring use cases (4 and 6)\n# is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Find-UseCase5Bounds($doc) {\n    $n = $doc.Paragraphs.Count\n    $startPara = 0\n    $endPara = 0\n    for ($i = 1; $i -le $n; $i++) {\n        $t = $doc.Paragraphs($i).Range.Text.Trim()\n        if ($t -eq \"Use case5\") {\n            $startPara = $i\n        } elseif ($startPara -gt 0 -and $i -gt $startPara -and $t -match \"^Use case\\d+$\") {\n            $endPara = $i\n            break\n        }\n    }\n    if ($startPara -eq 0) {\n        throw \"Could not locate the 'Use case 5' heading paragraph.\"\n    }\n    if ($endPara -eq 0) {\n        $endPara = $n + 1\n    }\n    return @{ Start = $startPara; End = $endPara }\n}\n\nfunction Get-ScopedRange($doc, $bounds) {\n    $s = $doc.Paragraphs($bounds.Start).Range.Start\n    if ($bounds.End -le $doc.Paragraphs.Count) {\n        $e = $doc.Paragraphs($bounds.End).Range.Start\n    } else {\n        $e = $doc.Content.End\n    }\n    return $doc.Range($s, $e)\n}\n\n# Finds (within the Use case 5 paragraph window) the first paragraph whose\n# trimmed text equals $ExpectedText and overwrites its text with\n# $NewText (the paragraph mark is preserved since we only touch the\n# paragraph's Range, which Word implicitly excludes the mark from when\n# setting .Text this way is fine because Range.Text assignment keeps the\n# paragraph mark intact).\nfunction Set-ParagraphText($doc, $bounds, [string]$ExpectedText, [string]$NewText) {\n    for ($i = $bounds.Start; $i -lt $bounds.End; $i++) {\n        $p = $doc.Paragraphs($i)\n        $t = $p.Range.Text.Trim()\n        if ($t -eq $ExpectedText) {\n            $r = $p.Range\n            $r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the paragraph mark\n            $r.Text = $NewText\n            return $true\n        }\n    }\n    return $false\n}\n\n$bounds = Find-UseCase5Bounds $d\n\n# 1) \"Adjunct uitleg\" -> \"Inverse uitleg\"\nSet-ParagraphText $d $bounds \"Adjunct uitleg\" \"Inverse uitleg\" | Out-Null\n\n# 2) & 4) Both \"Gebruiker weet hoe hij matrices moet transponeren\" bullets\n#    (Doelstelling + Postcondities) -> \"...inverteren\". Scope the\n#    replace-all to the Use case 5 paragraph range so the identical text\n#    in use case 4 / 6 is left alone.\n$scoped = Get-ScopedRange $d $bounds\n$scoped.Find.Execute(\"transponeren\", $false, $false, $false, $false, $false, $true, 1, $false, \"inverteren\", 2) | Out-Null\n\n# 3) \"Basiskennis matrices\" -> append \"+determinanten kunnen uitrekenen\"\n$bounds = Find-UseCase5Bounds $d\nfor ($i = $bounds.Start; $i -lt $bounds.End; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Basiskennis matrices\") {\n        $r = $p.Range\n        $r.MoveEnd(1, -1) | Out-Null\n        $r.Collapse(0) | Out-Null   # wdCollapseEnd = 0\n        $r.InsertAfter(\"+determinanten kunnen uitrekenen\")\n        break\n    }\n}\n\n# 5) \"Algemene definitie over transponeren van matrices\"\n#    -> \"Algemene definitie over inverteren van matrices\"\n$bounds = Find-UseCase5Bounds $d\nSet-ParagraphText $d $bounds \"Algemene definitie over transponeren van matrices\" \"Algemene definitie over inverteren van matrices\" | Out-Null\n\n# 6) \"Stap per stap animatie tonen\"\n#    -> \"Stap per stap bewerkingen uitleggen (adjunct berekenen)\"\n$bounds = Find-UseCase5Bounds $d\nSet-ParagraphText $d $bounds \"Stap per stap animatie tonen\" \"Stap per stap bewerkingen uitleggen (adjunct berekenen)\" | Out-Null\n\n# 7) Delete \"Gebruiker bekijkt de animatie en de extra uitleg stap per stap\"\n$bounds = Find-UseCase5Bounds $d\nfor ($i = $bounds.Start; $i -lt $bounds.End; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Gebruiker bekijkt de animatie en de extra uitleg stap per stap\") {\n        $p.Range.Delete() | Out-Null\n        break\n    }\n}\n\n# 8) Insert new paragraph \"Gebruiker klikt op voltooien\" right after\n#    \"Gebruiker klikt op next knop om naar volgende stap te gaan\" and\n#    before \"Gebruiker krijgt matrix\".\n$bounds = Find-UseCase5Bounds $d\nfor ($i = $bounds.Start; $i -lt $bounds.End; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Gebruiker klikt op next knop om naar volgende stap te gaan\") {\n        $p.Range.InsertParagraphAfter() | Out-Null\n        $d.Paragraphs($i + 1).Range.Text = \"Gebruiker klikt op voltooien\"\n        break\n    }\n}\n\n# 9) \"Geruiker stelt zelf de getransponeerde \" -> \"Geruiker stelt zelf de inverse op \"\n$bounds = Find-UseCase5Bounds $d\nSet-ParagraphText $d $bounds \"Geruiker stelt zelf de getransponeerde\" \"Geruiker stelt zelf de inverse op \" | Out-Null\n\n\"done\"\n"}
